# This script applies the "Updated capital structure database" edit:
#  - Laurentian Bank of Canada (row 5) is removed from the Canada / Banks (Regional) dataset
#  - The remaining three companies (rows 2-4) have their financial metrics recomputed
#    to reflect the new peer-group composition (columns D..AK)
#  - The now-unused debt_ebitda / net_debt_ebitda figures (columns AN, AP) are cleared
#  - The sheet dimension shrinks from A1:AQ5 to A1:AQ4 as a result

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 holds a text value ("2", was "3") - force text so it keeps its original string type
$ws.Range('B2').Value = "'2"

# --- Row 2: recomputed metrics ---
$ws.Range('D2').Value = 0.08385000000000001
$ws.Range('E2').Value = 0.07545
$ws.Range('I2').Value = 0
$ws.Range('J2').Value = 0
$ws.Range('K2').Value = 217.9
$ws.Range('L2').Value = 0.3372542950007739
$ws.Range('M2').Value = 80.7525
$ws.Range('N2').Value = 0.03838776383342841
$ws.Range('O2').Value = 0.3705943093162001
$ws.Range('P2').Value = 76.88249999999999
$ws.Range('Q2').Value = 0.03654806046776954
$ws.Range('R2').Value = 0.3528338687471317
$ws.Range('S2').Value = 3.870000000000005
$ws.Range('T2').Value = 0.04792421287266654
$ws.Range('U2').Value = 279.2
$ws.Range('V2').Value = 0.1327248526335805
$ws.Range('W2').Value = 0.09790290068218413
$ws.Range('X2').Value = 0.04501769685016251
$ws.Range('Y2').Value = 0.05288520383202162
$ws.Range('Z2').Value = 0.1593891849220446
$ws.Range('AA2').Value = 0
$ws.Range('AB2').Value = 0.03429712812302017
$ws.Range('AC2').Value = -0.03429712812302017
$ws.Range('AD2').Value = 1959.6
$ws.Range('AE2').Value = 0
$ws.Range('AF2').Value = 1959.6
$ws.Range('AG2').Value = 1680.4
$ws.Range('AH2').Value = 0.4822799763733018
$ws.Range('AI2').Value = 0.4208943683147901
$ws.Range('AJ2').Value = 0.4440803382663848
$ws.Range('AK2').Value = 0.3839510122012521

# --- Row 3: recomputed metrics ---
$ws.Range('D3').Value = 0.0677
$ws.Range('E3').Value = -0.0361
$ws.Range('I3').Value = 0
$ws.Range('J3').Value = 0
$ws.Range('K3').Value = 203.3
$ws.Range('L3').Value = 0.3359775243761362
$ws.Range('M3').Value = 79.17
$ws.Range('N3').Value = 0.0404671846248211
$ws.Range('O3').Value = 0.3894244958189867
$ws.Range('P3').Value = 75.3
$ws.Range('Q3').Value = 0.03848906154160703
$ws.Range('R3').Value = 0.3703885882931628
$ws.Range('S3').Value = 3.870000000000005
$ws.Range('T3').Value = 0.04888215233042825
$ws.Range('U3').Value = 85.59999999999999
$ws.Range('V3').Value = 0.04375383357186669
$ws.Range('W3').Value = 0.1046697214642434
$ws.Range('X3').Value = 0.0535631302426135
$ws.Range('Y3').Value = 0.05110659122162993
$ws.Range('Z3').Value = 0.1532791245535375
$ws.Range('AA3').Value = 0
$ws.Range('AB3').Value = 0.03388506737886177
$ws.Range('AC3').Value = -0.03388506737886177
$ws.Range('AD3').Value = 1947
$ws.Range('AE3').Value = 0
$ws.Range('AF3').Value = 1947
$ws.Range('AG3').Value = 1861.4
$ws.Range('AH3').Value = 0.4987959215043296
$ws.Range('AI3').Value = 0.4374003100217914
$ws.Range('AJ3').Value = 0.4875582796374875
$ws.Range('AK3').Value = 0.4263691962342808

# --- Row 4: recomputed metrics ---
$ws.Range('D4').Value = 0.1
$ws.Range('E4').Value = 0.187
$ws.Range('I4').Value = 0
$ws.Range('J4').Value = 0
$ws.Range('K4').Value = 14.6
$ws.Range('L4').Value = 0.3560975609756097
$ws.Range('M4').Value = 1.5825
$ws.Range('N4').Value = 0.01075067934782609
$ws.Range('O4').Value = 0.1083904109589041
$ws.Range('P4').Value = 1.5825
$ws.Range('Q4').Value = 0.01075067934782609
$ws.Range('R4').Value = 0.1083904109589041
$ws.Range('U4').Value = 193.6
$ws.Range('V4').Value = 1.315217391304348
$ws.Range('W4').Value = 0.09113607990012484
$ws.Range('X4').Value = 0.03647226345771152
$ws.Range('Y4').Value = 0.05466381644241333
$ws.Range('Z4').Value = 0.3871576959395657
$ws.Range('AA4').Value = 0
$ws.Range('AB4').Value = 0.03470918886717857
$ws.Range('AC4').Value = -0.03470918886717857
$ws.Range('AD4').Value = 12.6
$ws.Range('AE4').Value = 0
$ws.Range('AF4').Value = 12.6
$ws.Range('AG4').Value = -181
$ws.Range('AH4').Value = 0.07884856070087611
$ws.Range('AI4').Value = 0.06161369193154034
$ws.Range('AJ4').Value = 5.355029585798815
$ws.Range('AK4').Value = -16.60550458715596

# --- Clear debt_ebitda / net_debt_ebitda columns (AN, AP) for rows 2-4 ---
$ws.Range('AN2').ClearContents()
$ws.Range('AP2').ClearContents()
$ws.Range('AN3').ClearContents()
$ws.Range('AP3').ClearContents()
$ws.Range('AN4').ClearContents()
$ws.Range('AP4').ClearContents()

# --- Remove row 5 (Laurentian Bank of Canada) entirely ---
$ws.Rows(5).Delete()
